$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A8").PasteSpecial(-4122)

$ws.Range("A7").Value = "Smerovac"
$ws.Range("B7").Value = 2

$ws.Range("A8").Value = "Prepinac"
$ws.Range("B8").Value = 1

$ws.Range("E12").Select()
